# POCE fix in shelf placement
# The "num. of shelves max" / "num. ignored from bottom" values that used to
# live on the last data row (A5:D5 = 8,100,1,2) were a duplicate/POCE entry.
# Fix: fold the correct max-shelves (100) and ignored-from-bottom (2) values
# into row 4, then remove the now-redundant row 5 (shifting the rows below it
# up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 with the correct values that used to sit on row 5.
$ws.Range("B4").Value = 100
$ws.Range("D4").Value = 2

# Drop the now-duplicate row 5 entirely, shifting subsequent rows up.
$ws.Range("A5:D5").EntireRow.Delete()

# Leave the selection where the editor ended up.
$ws.Range("B10").Select()
